# Re-ran the CircaDB / CircadiPy cosinor analysis (sawtooth_0.5, cosinor_3)
# and refreshed the result table with the new simulation output.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -------- Row 2 --------
$ws.Range("B2").Value = [double]"1"
$ws.Range("E2").Value = [double]"25.11000000000049"
$ws.Range("G2").Value = [double]"0.001844347907995258"
$ws.Range("H2").Value = [double]"0.01043977479532889"

# p_reject (I2) no longer has a value for this row - blank it out like its
# neighbours (J2 / I3 / J3) that already carry no p_reject value.
$ws.Range("I2").ClearContents()

$ws.Range("K2").Value = [double]"4.496859984607298"
$ws.Range("L2").Value = "[1.306789425933781, 7.686930543280814]"
$ws.Range("M2").Value = [double]"0.005850815893534866"
$ws.Range("N2").Value = [double]"0.006032711322108364"
$ws.Range("O2").Value = [double]"-1.874263485024541"
$ws.Range("P2").Value = "[-2.6793162571156195, -1.0692107129334634]"
$ws.Range("Q2").Value = [double]"6.37658814595099e-06"
$ws.Range("R2").Value = [double]"6.37658814595099e-06"
$ws.Range("S2").Value = [double]"13.85883373332817"
$ws.Range("T2").Value = "[12.12632703759999, 15.591340429056347]"
$ws.Range("W2").Value = [double]"7.490270270270418"
$ws.Range("X2").Value = [double]"4.272972972973059"
$ws.Range("Y2").Value = [double]"10.70756756756778"

# -------- Row 3 --------
$ws.Range("E3").Value = [double]"24.68000000000042"
$ws.Range("G3").Value = [double]"0.00216392513673791"
$ws.Range("H3").Value = [double]"0.01043977479532889"
$ws.Range("K3").Value = [double]"4.568344022262792"
$ws.Range("L3").Value = "[1.3182890107369811, 7.818399033788603]"
$ws.Range("M3").Value = [double]"0.006032711322108364"
$ws.Range("N3").Value = [double]"0.006032711322108364"
$ws.Range("O3").Value = [double]"-2.830263651882697"
$ws.Range("P3").Value = "[-3.685632222229467, -1.9748950815359265]"
$ws.Range("Q3").Value = [double]"3.324707176233233e-10"
$ws.Range("R3").Value = [double]"6.649414352466465e-10"
$ws.Range("S3").Value = [double]"12.92350581852436"
$ws.Range("T3").Value = "[11.106216425952667, 14.740795211096057]"
$ws.Range("W3").Value = [double]"11.11711711711731"
$ws.Range("X3").Value = [double]"7.757277277277411"
$ws.Range("Y3").Value = [double]"14.4769569569572"
